$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "act.itt.p"
$ws.Range("B11").Value = 0.087
$ws.Range("E11").Value = 0.584
$ws.Range("I11").Value = "act.ppm.p"
$ws.Range("G11").Value = 0.145
$ws.Range("H11").Value = 0.6466
$ws.Range("F11").Value = "act.pp.p"
$ws.Range("J11").Value = 0.295
$ws.Range("K11").Value = 0.6759

$ws.Range("K12").Select()
